# Add Uncertainty to Model
# Update Q_cool (B3) with a new value reflecting the updated/uncertainty-adjusted model.
# B5 (T_cool) recalculates automatically since it holds the formula =B3/B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 88724.784780000002

$excel.Calculate()
